$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = -0.7091679257846845
$ws.Range("C2").Value = -0.2116494899511859
$ws.Range("D2").Value = -0.8069708324279797
$ws.Range("E2").Value = -0.442835019116573
$ws.Range("F2").Value = -0.02123990589338218
$ws.Range("G2").Value = -0.1707229344526847
$ws.Range("H2").Value = -0.1169158257196027
$ws.Range("B3").Value = -0.5825132856353434
$ws.Range("C3").Value = -1.177834628112137
$ws.Range("D3").Value = -0.8136988148007305
$ws.Range("E3").Value = -0.3921037015775397
$ws.Range("F3").Value = -0.5415867301368422
$ws.Range("G3").Value = -0.4877796214037602
$ws.Range("B4").Value = -0.6994291822578886
$ws.Range("C4").Value = -0.3352933689464819
$ws.Range("D4").Value = 0.08630174427670884
$ws.Range("E4").Value = -0.06318128428259367
$ws.Range("F4").Value = -0.009374175549511699
$ws.Range("G4").Value = -0.334123210679779
$ws.Range("H4").Value = 0.001243794298757539
$ws.Range("I4").Value = 0.09478491567901273
$ws.Range("J4").Value = -0.3880619836346724
$ws.Range("B5").Value = -0.1970460893006987
$ws.Range("C5").Value = 0.2245490239224921
$ws.Range("D5").Value = 0.07506599536318959
$ws.Range("E5").Value = 0.1288731040962716
$ws.Range("F5").Value = -0.1958759310339957
$ws.Range("G5").Value = 0.1394910739445408
$ws.Range("H5").Value = 0.233032195324796
$ws.Range("I5").Value = -0.2498147039888892
$ws.Range("B6").Value = 0.1735200387801499
$ws.Range("C6").Value = 0.02403701022084742
$ws.Range("D6").Value = 0.0778441189539294
$ws.Range("E6").Value = -0.2469049161763379
$ws.Range("F6").Value = 0.08846208880219864
$ws.Range("G6").Value = 0.1820032101824538
$ws.Range("H6").Value = -0.3008436891312313
$ws.Range("B7").Value = -0.009228378467612708
$ws.Range("C7").Value = 0.04457873026546927
$ws.Range("D7").Value = -0.280170304864798
$ws.Range("E7").Value = 0.05519670011373851
$ws.Range("F7").Value = 0.1487378214939937
$ws.Range("G7").Value = -0.3341090778196915
$ws.Range("B8").Value = 0.1858859418247861
$ws.Range("C8").Value = -0.1388630933054812
$ws.Range("D8").Value = 0.1965039116730553
$ws.Range("E8").Value = 0.2900450330533105
$ws.Range("F8").Value = -0.1928018662603747
$ws.Range("G8").Value = -0.08321328366289282
$ws.Range("H8").Value = -0.05405384115675094
$ws.Range("I8").Value = -0.2527483093463515
$ws.Range("B9").Value = -0.2674987412282675
$ws.Range("C9").Value = 0.067868263750269
$ws.Range("D9").Value = 0.1614093851305242
$ws.Range("E9").Value = -0.321437514183161
$ws.Range("F9").Value = -0.2118489315856791
$ws.Range("G9").Value = -0.1826894890795372
$ws.Range("H9").Value = -0.3813839572691378
$ws.Range("B10").Value = 0.07761257013415129
$ws.Range("C10").Value = 0.1711536915144065
$ws.Range("D10").Value = -0.3116932077992787
$ws.Range("E10").Value = -0.2021046252017968
$ws.Range("F10").Value = -0.172945182695655
$ws.Range("G10").Value = -0.3716396508852555
$ws.Range("B11").Value = -0.1119472883713244
$ws.Range("C11").Value = -0.5947941876850096
$ws.Range("D11").Value = -0.4852056050875277
$ws.Range("E11").Value = -0.4560461625813858
$ws.Range("F11").Value = -0.6547406307709864
$ws.Range("B12").Value = -0.2799851348954449
$ws.Range("C12").Value = -0.170396552297963
$ws.Range("D12").Value = -0.1412371097918211
$ws.Range("E12").Value = -0.3399315779814217
$ws.Range("B13").Value = -0.1408857896223169
$ws.Range("C13").Value = -0.111726347116175
$ws.Range("D13").Value = -0.3104208153057756
$ws.Range("B14").Value = -0.2487487154747626
$ws.Range("C14").Value = -0.4474431836643632
$ws.Range("B15").Value = -0.3133716091653966
